{"js": "const replacements = [\n  [\"2025-04-06 Sunday\", \"2025-04-07 Monday\"],\n  [\"90\u00d775=6750\", \"64\u00d786=5504\"],\n  [\"14\u00d771=994\", \"64\u00d736=2304\"],\n  [\"21\u00d768=1428\", \"32\u00d771=2272\"],\n  [\"92\u00d778=7176\", \"12\u00d723=276\"],\n  [\"14\u00d798=1372\", \"90\u00d715=1350\"],\n  [\"73\u00d786=6278\", \"31\u00d793=2883\"],\n  [\"92\u00d733=3036\", \"29\u00d785=2465\"],\n  [\"31\u00d752=1612\", \"36\u00d712=432\"],\n  [\"43\u00d770=3010\", \"13\u00d749=637\"],\n  [\"33\u00d745=1485\", \"51\u00d731=1581\"],\n  [\"91\u00d743=3913\", \"15\u00d720=300\"],\n  [\"17\u00d711=187\", \"37\u00d768=2516\"],\n  [\"54\u00d761=3294\", \"46\u00d761=2806\"],\n  [\"57\u00d785=4845\", \"32\u00d727=864\"],\n  [\"74\u00d761=4514\", \"54\u00d746=2484\"],\n  [\"57\u00d749=2793\", \"98\u00d747=4606\"],\n  [\"38\u00d786=3268\", \"20\u00d798=1960\"],\n  [\"42\u00d742=1764\", \"83\u00d793=7719\"],\n  [\"35\u00d724=840\", \"79\u00d771=5609\"],\n  [\"41\u00d736=1476\", \"89\u00d717=1513\"],\n  [\"39\u00d767=2613\", \"41\u00d733=1353\"],\n  [\"93\u00d776=7068\", \"98\u00d721=2058\"],\n  [\"99\u00d797=9603\", \"75\u00d716=1200\"],\n  [\"89\u00d772=6408\", \"34\u00d746=1564\"],\n  [\"69\u00d783=5727\", \"53\u00d773=3869\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @{old=\"2025-04-06 Sunday\"; new=\"2025-04-07 Monday\"},\n    @{old=\"90\u00d775=6750\"; new=\"64\u00d786=5504\"},\n    @{old=\"14\u00d771=994\"; new=\"64\u00d736=2304\"},\n    @{old=\"21\u00d768=1428\"; new=\"32\u00d771=2272\"},\n    @{old=\"92\u00d778=7176\"; new=\"12\u00d723=276\"},\n    @{old=\"14\u00d798=1372\"; new=\"90\u00d715=1350\"},\n    @{old=\"73\u00d786=6278\"; new=\"31\u00d793=2883\"},\n    @{old=\"92\u00d733=3036\"; new=\"29\u00d785=2465\"},\n    @{old=\"31\u00d752=1612\"; new=\"36\u00d712=432\"},\n    @{old=\"43\u00d770=3010\"; new=\"13\u00d749=637\"},\n    @{old=\"33\u00d745=1485\"; new=\"51\u00d731=1581\"},\n    @{old=\"91\u00d743=3913\"; new=\"15\u00d720=300\"},\n    @{old=\"17\u00d711=187\"; new=\"37\u00d768=2516\"},\n    @{old=\"54\u00d761=3294\"; new=\"46\u00d761=2806\"},\n    @{old=\"57\u00d785=4845\"; new=\"32\u00d727=864\"},\n    @{old=\"74\u00d761=4514\"; new=\"54\u00d746=2484\"},\n    @{old=\"57\u00d749=2793\"; new=\"98\u00d747=4606\"},\n    @{old=\"38\u00d786=3268\"; new=\"20\u00d798=1960\"},\n    @{old=\"42\u00d742=1764\"; new=\"83\u00d793=7719\"},\n    @{old=\"35\u00d724=840\"; new=\"79\u00d771=5609\"},\n    @{old=\"41\u00d736=1476\"; new=\"89\u00d717=1513\"},\n    @{old=\"39\u00d767=2613\"; new=\"41\u00d733=1353\"},\n    @{old=\"93\u00d776=7068\"; new=\"98\u00d721=2058\"},\n    @{old=\"99\u00d797=9603\"; new=\"75\u00d716=1200\"},\n    @{old=\"89\u00d772=6408\"; new=\"34\u00d746=1564\"},\n    @{old=\"69\u00d783=5727\"; new=\"53\u00d773=3869\"}\n)\n\nforeach ($p in $pairs) {\n    $rng = $d.Content\n    $rng.Find.Execute($p.old, $false, $false, $false, $false, $false, $true, 1, $false, $p.new, 2)\n}\n"}
